# Swap the taxon-observation data between row 12 and row 13, while keeping
# the location/date/administrative columns (P, S, T, U, V, W, Y, AA, AD, AE,
# AG, AT, AW, AY, I) tied to their original row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture current ("before") values of the columns that move between rows ---
$A12 = $ws.Range("A12").Value2
$B12 = $ws.Range("B12").Value2
$D12 = $ws.Range("D12").Value2
$E12 = $ws.Range("E12").Value2
$F12 = $ws.Range("F12").Value2
$G12 = $ws.Range("G12").Value2
$H12 = $ws.Range("H12").Value2
$K12 = $ws.Range("K12").Value2
$L12 = $ws.Range("L12").Value2
$M12 = $ws.Range("M12").Value2
$N12 = $ws.Range("N12").Value2
$Q12 = $ws.Range("Q12").Value2
$R12 = $ws.Range("R12").Value2
$AX12 = $ws.Range("AX12").Value2

$A13 = $ws.Range("A13").Value2
$B13 = $ws.Range("B13").Value2
$D13 = $ws.Range("D13").Value2
$E13 = $ws.Range("E13").Value2
$F13 = $ws.Range("F13").Value2
$G13 = $ws.Range("G13").Value2
$H13 = $ws.Range("H13").Value2
$K13 = $ws.Range("K13").Value2
$L13 = $ws.Range("L13").Value2
$M13 = $ws.Range("M13").Value2
$N13 = $ws.Range("N13").Value2
$Q13 = $ws.Range("Q13").Value2
$R13 = $ws.Range("R13").Value2
$AX13 = $ws.Range("AX13").Value2

# --- Write row 12 with row 13's former data ---
$ws.Range("A12").Value = $A13
$ws.Range("B12").Value = $B13
$ws.Range("D12").Value = $D13
$ws.Range("E12").Value = $E13
$ws.Range("F12").Value = $F13
$ws.Range("G12").Value = $G13
$ws.Range("H12").Value = $H13
$ws.Range("K12").Value = $K13
$ws.Range("L12").Value = $L13
$ws.Range("M12").Value = $M13
$ws.Range("N12").Value = $N13
$ws.Range("Q12").Value = $Q13
$ws.Range("R12").Value = $R13
$ws.Range("AX12").Value = $AX13

# --- Write row 13 with row 12's former data ---
$ws.Range("A13").Value = $A12
$ws.Range("B13").Value = $B12
$ws.Range("D13").Value = $D12
$ws.Range("E13").Value = $E12
$ws.Range("F13").Value = $F12
$ws.Range("G13").Value = $G12
$ws.Range("H13").Value = $H12
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("Q13").Value = $Q12
$ws.Range("R13").Value = $R12
$ws.Range("AX13").Value = $AX12
